$d = $word.ActiveDocument

# --- Create the new "centre_justified" paragraph style (custom style,
#     styleId "centrejustified") and its linked "centre_justified Char"
#     character style (styleId "centrejustifiedChar"), modelled on
#     Heading4 / Heading4Char but centre aligned. ---

$paraStyle = $d.Styles.Add("centrejustified", 1)
$paraStyle.NameLocal = "centre_justified"
$paraStyle.QuickStyle = $true
$paraStyle.ParagraphFormat.Alignment = 1
$paraStyle.Font.Name = "Arial"
$paraStyle.Font.Bold = $true
$paraStyle.Font.BoldBi = $true
$paraStyle.Font.ItalicBi = $true
$paraStyle.Font.Spacing = -1
$paraStyle.Font.Size = 12
$paraStyle.Font.SizeBi = 11
$paraStyle.Font.LanguageIDFarEast = "en-US"

$charStyle = $d.Styles.Add("centrejustifiedChar", 2)
$charStyle.NameLocal = "centre_justified Char"
$charStyle.BaseStyle = $d.Styles("Heading4Char")
$charStyle.Font.Bold = $true
$charStyle.Font.BoldBi = $true
$charStyle.Font.ItalicBi = $true

$paraStyle.LinkStyle = $charStyle
$charStyle.LinkStyle = $paraStyle

# --- Apply the new style to the first (empty) paragraph, replacing
#     "Heading4". ---

$d.Paragraphs(1).Style = $paraStyle
